$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Accent1 (theme color 4 / "4F81BD") as a plain RGB OLE color - used where
# ThemeColor assignment isn't honoured by this host (Borders).
$accent1Rgb = 12419407   # 0xBD814F => R=4F G=81 B=BD

# ---------------------------------------------------------------------
# 1. Rename sheet
# ---------------------------------------------------------------------
$ws.Name = "Sheet1"

# ---------------------------------------------------------------------
# 2. Drop the now-unused data rows (rows 5-10); keep the 3 header/template
#    rows and a single sample data row (row 4)
# ---------------------------------------------------------------------
$ws.Rows("5:10").Delete() | Out-Null

# ---------------------------------------------------------------------
# 3. Resize the table to the new header/data block BEFORE the header text
#    changes, so its column names track row 3 once that row gets its new
#    field-name text below.
# ---------------------------------------------------------------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A3:D4")) | Out-Null

# ---------------------------------------------------------------------
# 4. Copy the pre-existing "Chinese label" look (bold header font + accent
#    fill, originally on row 3) up onto row 1, and stamp the "type name"
#    look (italic font + accent fill, originally on row 2) across column D
#    too (columns A-C already carry it).
# ---------------------------------------------------------------------
$ws.Range("A3:C3").Copy() | Out-Null
$ws.Range("A1:D1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("A2:C2").Copy() | Out-Null
$ws.Range("D2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 5. Rewrite the cell text for all four rows.
# ---------------------------------------------------------------------
# Row 1: Chinese field labels
$ws.Range("A1").Value = "序列"
$ws.Range("B1").Value = "英文名"
$ws.Range("C1").Value = "道具列表"
$ws.Range("D1").Value = "货币类型"

# Row 2: type names
$ws.Range("A2").Value = "int"
$ws.Range("B2").Value = "string"
$ws.Range("C2").Value = "string[]"
$ws.Range("D2").Value = "int"

# Row 3: English field names (table header row)
$ws.Range("A3").Value = "Id"
$ws.Range("B3").Value = "Ename"
$ws.Range("C3").Value = "SellTable"
$ws.Range("D3").Value = "MoneyType"

# Row 4: sample data row
$ws.Range("A4").Value = 44000001
$ws.Range("B4").Value = "sdaolai"
$ws.Range("C4").Value = "zzwandou;zzyumi;zzpingguo;zzlanmei"
$ws.Range("D4").Value = 0

# ---------------------------------------------------------------------
# 6. Rename & restyle the table
# ---------------------------------------------------------------------
$lo.Name = "表3"
$lo.TableStyle = "TableStyleLight13"

# ---------------------------------------------------------------------
# 7. Borders: rows 1 & 2 get a thin accent-coloured top border running the
#    full width, with the left edge also boxed on column A and the right
#    edge boxed on column D.
# ---------------------------------------------------------------------
foreach ($r in 1, 2) {
    $left = $ws.Range("A$r")
    $left.Borders.Item(7).Color = $accent1Rgb   # xlEdgeLeft
    $left.Borders.Item(7).LineStyle = 1
    $left.Borders.Item(8).Color = $accent1Rgb   # xlEdgeTop
    $left.Borders.Item(8).LineStyle = 1

    $mid = $ws.Range("B$r:C$r")
    $mid.Borders.Item(8).Color = $accent1Rgb
    $mid.Borders.Item(8).LineStyle = 1

    $right = $ws.Range("D$r")
    $right.Borders.Item(10).Color = $accent1Rgb  # xlEdgeRight
    $right.Borders.Item(10).LineStyle = 1
    $right.Borders.Item(8).Color = $accent1Rgb
    $right.Borders.Item(8).LineStyle = 1
}

# Row 4 (data row) - plain thin top border across the whole row
$ws.Range("A4:D4").Borders.Item(8).Color = $accent1Rgb
$ws.Range("A4:D4").Borders.Item(8).LineStyle = 1

# ---------------------------------------------------------------------
# 8. Row 1: vertical "stacked" text + accent fill on top of the copied font
# ---------------------------------------------------------------------
$ws.Range("A1:D1").Orientation = -4166   # xlVertical -> textRotation=255

# ---------------------------------------------------------------------
# 9. Row 3 (table header): bold white-on-accent look
# ---------------------------------------------------------------------
$hdr = $ws.Range("A3:D3")
$hdr.Font.ThemeColor = 2   # xlThemeColorLight1 -> theme 0 (white)
$hdr.Font.Bold = $true
$hdr.Interior.ThemeColor = 5   # xlThemeColorAccent1 -> theme 4

Write-Host "formatting done"

# ---------------------------------------------------------------------
# 10. Column widths / row height / selection
# ---------------------------------------------------------------------
$ws.Columns("C:D").ColumnWidth = 12.16

$ws.Rows("1").RowHeight = 55.5

$ws.Range("C4").Select() | Out-Null

Write-Host "layout done"
